$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "_validation_data" sheet: the "mass_analysis_polarity" value column (C)
#    is removed entirely; ms_scan_mode / preparation_matrix / duration-unit
#    values shift one column to the left (D->C, E->D, F->E).
# ---------------------------------------------------------------------------
$wsValidation = $wb.Worksheets.Item("_validation_data")
$wsValidation.Columns.Item(3).Delete()

# ---------------------------------------------------------------------------
# 2. "Non-Standard Value" sheet: the "mass_analysis_polarity" row (row 4) is
#    removed; the "ms_scan_mode" and "preparation_matrix" rows shift up.
# ---------------------------------------------------------------------------
$wsNonStandard = $wb.Worksheets.Item("Non-Standard Value")
$wsNonStandard.Rows.Item(4).Delete()

# Re-apply the AutoFilter so the ref shrinks from A1:E6 to A1:E5.
$wsNonStandard.AutoFilterMode = $false
[void]$wsNonStandard.Range("A1:E5").AutoFilter()

# Fix up the data-validation source ranges on "Non-Standard Value" so they
# keep pointing at the right columns of "_validation_data" after the column
# deletion above (ms_scan_mode: C, preparation_matrix: D).
$wsNonStandard.Range("D4").Validation.Formula1 = "_validation_data!`$C`$1:`$C`$3"
$wsNonStandard.Range("D5").Validation.Formula1 = "_validation_data!`$D`$1:`$D`$8"

# ---------------------------------------------------------------------------
# 3. "Missing Required Value" sheet: fix up the data-validation source
#    ranges for preparation_matrix (E->D) and the storage-duration unit
#    (F->E) to match the shifted "_validation_data" columns.
# ---------------------------------------------------------------------------
$wsMissing = $wb.Worksheets.Item("Missing Required Value")
$wsMissing.Range("C3").Validation.Formula1 = "_validation_data!`$D`$1:`$D`$8"
$wsMissing.Range("C4").Validation.Formula1 = "_validation_data!`$E`$1:`$E`$5"

# ---------------------------------------------------------------------------
# 4. Workbook-level hidden _FilterDatabase defined name for "Non-Standard
#    Value" needs to shrink the same way the AutoFilter ref did.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "Non-Standard Value!_FilterDatabase") {
        $n.RefersTo = "='Non-Standard Value'!`$A`$1:`$E`$5"
    }
}
